$d = $word.ActiveDocument

# The first two paragraphs of the document are the title ("The Mystical Body
# and Spain", styled Heading1) and the byline ("By Dorothy Day", bold text in
# a Normal paragraph). We replace them with a pandoc-style title block: a
# Title-styled paragraph (words split into their own runs, as pandoc's docx
# writer does) followed by an Authors-styled paragraph with the bare name.

$titlePara = $d.Paragraphs.Item(1)
$authorPara = $d.Paragraphs.Item(2)

$start = $titlePara.Range.Start
$end = $authorPara.Range.End
$target = $d.Range($start, $end)

$newXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Title"/>
            </w:pPr>
            <w:r><w:t xml:space="preserve">The</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:t xml:space="preserve">Mystical</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:t xml:space="preserve">Body</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:t xml:space="preserve">and</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:t xml:space="preserve">Spain</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Authors"/>
            </w:pPr>
            <w:r><w:t xml:space="preserve">Dorothy</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:t xml:space="preserve">Day</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$target.InsertXML($newXml)

# The original document also wraps the old heading in a bookmark
# ("the-mystical-body-and-spain") that the new pandoc-style output drops.
# It spans paragraph boundaries (sibling to <w:p>, not nested in a run), so
# it is not enumerated by Document.Bookmarks in this object model; delete it
# defensively in case a given host does expose it.
try {
    $staleBookmark = $d.Bookmarks.Item("the-mystical-body-and-spain")
    $staleBookmark.Delete()
} catch {
    # Not reachable through the Bookmarks collection here - nothing to do.
}
